{"js": "// 1) Insert a new \"Last run on: 2020-12-15 10:09:20\" Subtitle paragraph\n//    right after the Title paragraph (\"R Markdown Output\").\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst titlePara = paragraphs.items[0];\nconst lastRunPara = titlePara.insertParagraph(\n  \"Last run on: 2020-12-15 10:09:20\",\n  \"After\"\n);\nlastRunPara.style = \"Subtitle\";\nawait context.sync();\n\n// 2) Remove the trailing \"Back to Gallery\" section: the Heading1 link\n//    paragraph, the SourceCode `include_url(...)` snippet paragraph, and\n//    the paragraph holding the linked gallery image - the last three\n//    paragraphs of the document, right before the final section break.\nconst allParagraphs = context.document.body.paragraphs;\nallParagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = allParagraphs.items;\nconst count = items.length;\nfor (let i = count - 1; i >= count - 3; i--) {\n  items[i].delete();\n}\nawait context.sync();\n", "ps1": "# 1) Insert a new \"Last run on: 2020-12-15 10:09:20\" Subtitle paragraph\n#    right after the Title paragraph (\"R Markdown Output\").\n$d = $word.ActiveDocument\n\n$titlePara = $d.Paragraphs.Item(1)\n$afterTitle = $titlePara.Range\n$afterTitle.Collapse(0)  # wdCollapseEnd\n$afterTitle.InsertParagraphAfter()\n\n$newPara = $d.Paragraphs.Item(2)\n$newPara.Range.Text = \"Last run on: 2020-12-15 10:09:20\"\n$newPara.Style = \"Subtitle\"\n\n# 2) Remove the trailing \"Back to Gallery\" section: the Heading1 link\n#    paragraph, the SourceCode `include_url(...)` snippet paragraph, and\n#    the paragraph holding the linked gallery image - the last three\n#    paragraphs of the document, right before the final section break.\n$count = $d.Paragraphs.Count\n$startPara = $d.Paragraphs.Item($count - 2)\n$startOfSection = $startPara.Range.Start\n$endOfDoc = $d.Content.End\n$deleteRange = $d.Range($startOfSection, $endOfDoc)\n$deleteRange.Delete()\n"}
